$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values are being refreshed with a new day's electricity spot prices.
$ws.Range("A2").Value = 45935
$ws.Range("B2").Value = 35.02
$ws.Range("C2").Value = 33.51
$ws.Range("D2").Value = 29.05
$ws.Range("E2").Value = 28.93
$ws.Range("F2").Value = 28.2
$ws.Range("G2").Value = 28.2
$ws.Range("H2").Value = 29.19
$ws.Range("I2").Value = 28.1
$ws.Range("J2").Value = 24.85
$ws.Range("K2").Value = 5.27
$ws.Range("L2").Value = -0.01
$ws.Range("M2").Value = -0.64
$ws.Range("N2").Value = -1.12
$ws.Range("O2").Value = -2.22
$ws.Range("P2").Value = -2.45
$ws.Range("Q2").Value = -1.21
$ws.Range("R2").Value = -0.63
$ws.Range("S2").Value = -0.12
$ws.Range("T2").Value = 13.87
$ws.Range("U2").Value = 52.43
$ws.Range("V2").Value = 70.34
$ws.Range("W2").Value = 83.13
$ws.Range("X2").Value = 60.49
$ws.Range("Y2").Value = 41.13
$ws.Range("Z2").Value = 24.3
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 63.77
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 76.74
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 50.81
$ws.Range("AG2").Value = "9h-18h"
